$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Baza podataka" -> "Analiza sadržaja"
$ws.Range("B4").Value = "Analiza sadržaja"

# Move selection to the edited cell
$ws.Range("B4").Select()
